$d = $word.ActiveDocument

# --- 1. Remove the trailing duplicate bold title paragraph near the end of the doc ---
# (Do this before inserting the new "Meta description" paragraph up top so the
#  title text is still unique at this point - only the true duplicate carries
#  the "Normal" style, the real title keeps "Heading 1".)
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Play Crime City Online for Free - Slot Review" -and $p.Style.NameLocal -ne "Heading 1") {
        $p.Range.Delete()
    }
}

# --- 2. Replace the trailing meta-description text with the image prompt text ---
$d.Content.Find.Execute("Read our review of Crime City, a slot machine with a cyberpunk theme, 243 ways to win, and flexible betting options. Play for free now!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Prompt: Design a cartoon-style feature image for Crime City featuring a Maya warrior with glasses. The feature image should capture the theme of Crime City, which is crime in a futuristic metropolis. The background of the image should represent a dark city skyline with glaring neon lights. The Maya warrior should be shown in a happy and confident pose, holding a badge or a gun to represent the symbols in the game. The warrior should be depicted wearing glasses, which is a unique detail that sets the image apart from typical crime-themed images. The glasses should contribute to the cool and edgy vibe of the image. The color scheme should be mainly black with contrasting neon colors to highlight the warrior and the symbols. The image should be eye-catching and convey the excitement and thrill of playing Crime City.",
    2)

# --- 3. Insert a new "Meta description" paragraph right after the title paragraph ---
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
# Collapse to just before the paragraph's trailing paragraph mark so the inserted
# content becomes its own clean paragraph (no inherited pStyle / rsid noise).
$insertPos = $titleRange.End - 1
$insertRng = $d.Range($insertPos, $insertPos)

$metaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Crime City, a slot machine with a cyberpunk theme, 243 ways to win, and flexible betting options. Play for free now!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertRng.InsertXML($metaXml)
